# Update NATMI ligand/receptor (Il7-Il2rg) TPM results on Sheet1.
# - Re-computes the existing ECs/FAPs/MuSCs sending-cluster rows (2-13) with the
#   new TPM-based statistics.
# - Adds a 4th sending cluster, "Resolving-Mac" (rows 14-17), paired against all
#   4 target clusters, growing the sheet from A1:T13 to A1:T17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Il7","Il2rg","ECs",3,1,0.5282253333333333,1.584676,0.5375365413017625,0.5375365413017625,3,1,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,23.07624383387111,207.68619450484,0.250275411269706,0.250275411269706),
  @("ECs","Il7","Il2rg","FAPs",3,1,0.5282253333333333,1.584676,0.5375365413017625,0.5375365413017625,3,1,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,0.6366196367848889,5.729576731063999,0.006904513688005211,0.006904513688005211),
  @("ECs","Il7","Il2rg","MuSCs",3,1,0.5282253333333333,1.584676,0.5375365413017625,0.5375365413017625,3,1,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,2.631208883883111,23.680879954948,0.02853700499488438,0.02853700499488438),
  @("ECs","Il7","Il2rg","Resolving-Mac",3,1,0.5282253333333333,1.584676,0.5375365413017625,0.5375365413017625,3,1,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,23.21862433134444,208.9676189821,0.251819611349167,0.251819611349167),
  @("FAPs","Il7","Il2rg","ECs",2,0.6666666666666666,0.2859449999999999,0.8578349999999999,0.2909854499642813,0.2909854499642813,3,1,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,12.49189716335,112.42707447015,0.1354819580952499,0.1354819580952499),
  @("FAPs","Il7","Il2rg","FAPs",2,0.6666666666666666,0.2859449999999999,0.8578349999999999,0.2909854499642813,0.2909854499642813,3,1,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,0.3446222484099999,3.101600235689999,0.003737630594235004,0.003737630594235005),
  @("FAPs","Il7","Il2rg","MuSCs",2,0.6666666666666666,0.2859449999999999,0.8578349999999999,0.2909854499642813,0.2909854499642813,3,1,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,1.424356191995,12.819205727955,0.01544797906940387,0.01544797906940387),
  @("FAPs","Il7","Il2rg","Resolving-Mac",2,0.6666666666666666,0.2859449999999999,0.8578349999999999,0.2909854499642813,0.2909854499642813,3,1,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,12.568972208375,113.120749875375,0.1363178822053925,0.1363178822053926),
  @("MuSCs","Il7","Il2rg","ECs",1,0.3333333333333333,0.1276103333333333,0.382831,0.1298597641682559,0.1298597641682559,3,1,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,5.574831387087777,50.17348248379,0.06046231909348839,0.0604623190934884),
  @("MuSCs","Il7","Il2rg","FAPs",1,0.3333333333333333,0.1276103333333333,0.382831,0.1298597641682559,0.1298597641682559,3,1,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,0.1537965692482222,1.384169123234,0.001668014079655856,0.001668014079655856),
  @("MuSCs","Il7","Il2rg","MuSCs",1,0.3333333333333333,0.1276103333333333,0.382831,0.1298597641682559,0.1298597641682559,3,1,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,0.6356556975847778,5.720901278263,0.006894059201500232,0.006894059201500232),
  @("MuSCs","Il7","Il2rg","Resolving-Mac",1,0.3333333333333333,0.1276103333333333,0.382831,0.1298597641682559,0.1298597641682559,3,1,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,5.609228114386111,50.483053029475,0.0608353717936114,0.0608353717936114),
  @("Resolving-Mac","Il7","Il2rg","ECs",1,0.3333333333333333,0.04089733333333333,0.122692,0.0416182445657004,0.0416182445657004,3,1,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,1.786655763364444,16.07990187028,0.01937733060859303,0.01937733060859303),
  @("Resolving-Mac","Il7","Il2rg","FAPs",1,0.3333333333333333,0.04089733333333333,0.122692,0.0416182445657004,0.0416182445657004,3,1,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,0.04928965698755556,0.443606912888,0.0005345752654856484,0.0005345752654856484),
  @("Resolving-Mac","Il7","Il2rg","MuSCs",1,0.3333333333333333,0.04089733333333333,0.122692,0.0416182445657004,0.0416182445657004,3,1,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,0.2037187919684444,1.833469127716,0.002209449891859506,0.002209449891859506),
  @("Resolving-Mac","Il7","Il2rg","Resolving-Mac",1,0.3333333333333333,0.04089733333333333,0.122692,0.0416182445657004,0.0416182445657004,3,1,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,1.797679435077778,16.1791149157,0.01949688879976222,0.01949688879976222)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}
